$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 53-54 (previously the last two rows) and append
# --- four brand-new rows (55-58) describing newly generated rooms.
# Row 53
$ws.Cells.Item(53, 1).Value = 2052
$ws.Cells.Item(53, 2).Value = 1200
$ws.Cells.Item(53, 3).Value = "Superior"
$ws.Cells.Item(53, 4).Value = "Normal"
$ws.Cells.Item(53, 5).Value = "A"
$ws.Cells.Item(53, 6).Value = 1
$ws.Cells.Item(53, 7).Value = 1
$ws.Cells.Item(53, 8).Value = 1200
$ws.Cells.Item(53, 9).Value = "Unoccupied"
$ws.Cells.Item(53, 10).Value = "24-04-2020 03:28:54"

# Row 54
$ws.Cells.Item(54, 1).Value = 2053
$ws.Cells.Item(54, 2).Value = 1200
$ws.Cells.Item(54, 3).Value = "Superior"
$ws.Cells.Item(54, 4).Value = "Normal"
$ws.Cells.Item(54, 5).Value = "A"
$ws.Cells.Item(54, 6).Value = 1
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 1200
$ws.Cells.Item(54, 9).Value = "Unoccupied"
$ws.Cells.Item(54, 10).Value = "24-04-2020 03:30:52"

# Row 55 (new)
$ws.Cells.Item(55, 1).Value = 2054
$ws.Cells.Item(55, 2).Value = 1300
$ws.Cells.Item(55, 3).Value = "Superior"
$ws.Cells.Item(55, 4).Value = "Normal"
$ws.Cells.Item(55, 5).Value = "A"
$ws.Cells.Item(55, 6).Value = 1
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 1300
$ws.Cells.Item(55, 9).Value = "Unoccupied"
$ws.Cells.Item(55, 10).Value = "24-04-2020 03:31:01"

# Row 56 (new)
$ws.Cells.Item(56, 1).Value = 2055
$ws.Cells.Item(56, 2).Value = 1102
$ws.Cells.Item(56, 3).Value = "Superior"
$ws.Cells.Item(56, 4).Value = "Normal"
$ws.Cells.Item(56, 5).Value = "A"
$ws.Cells.Item(56, 6).Value = 1
$ws.Cells.Item(56, 7).Value = 1
$ws.Cells.Item(56, 8).Value = 1102
$ws.Cells.Item(56, 9).Value = "Unoccupied"
$ws.Cells.Item(56, 10).Value = "24-04-2020 03:33:59"

# Row 57 (new)
$ws.Cells.Item(57, 1).Value = 2056
$ws.Cells.Item(57, 2).Value = 9999
$ws.Cells.Item(57, 3).Value = "Superior"
$ws.Cells.Item(57, 4).Value = "Normal"
$ws.Cells.Item(57, 5).Value = "A"
$ws.Cells.Item(57, 6).Value = 1
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 9999
$ws.Cells.Item(57, 9).Value = "Unoccupied"
$ws.Cells.Item(57, 10).Value = "24-04-2020 03:34:24"

# Row 58 (new)
$ws.Cells.Item(58, 1).Value = 2057
$ws.Cells.Item(58, 2).Value = 9999
$ws.Cells.Item(58, 3).Value = "Superior"
$ws.Cells.Item(58, 4).Value = "Normal"
$ws.Cells.Item(58, 5).Value = "A"
$ws.Cells.Item(58, 6).Value = 1
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 9999
$ws.Cells.Item(58, 9).Value = "Unoccupied"
$ws.Cells.Item(58, 10).Value = "24-04-2020 03:35:03"

# --- The header cell I1 ("Status") gets overwritten with "Occupied" last,
# --- matching the order new shared strings were minted in the source edit.
$ws.Cells.Item(1, 9).Value = "Occupied"
